$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 713.5714
$ws.Range("I111").Value = 373.75
$ws.Range("K111").Value = 1121.25
$ws.Range("M111").Value = 1945.75
$ws.Range("H132").Value = 1635.7805
$ws.Range("I132").Value = 1212.5264
$ws.Range("K132").Value = 3637.5792
$ws.Range("M132").Value = -1107.5792
$ws.Range("H133").Value = 73297.586
$ws.Range("J133").Value = 73297.586
$ws.Range("L133").Value = 73297.586
$ws.Range("N133").Value = -83417.586
$ws.Range("H134").Value = 54997.777
$ws.Range("J134").Value = 54997.777
$ws.Range("L134").Value = 54997.777
$ws.Range("N134").Value = -65137.777
$ws.Range("H139").Value = 98406
$ws.Range("J139").Value = 98406
$ws.Range("L139").Value = 98406
$ws.Range("N139").Value = -108686
$ws.Range("H140").Value = 80776.664
$ws.Range("J140").Value = 80776.664
$ws.Range("L140").Value = 80776.664
$ws.Range("N140").Value = -91136.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 23266.666
$ws.Range("J130").Value = 23266.666
$ws.Range("L130").Value = 23266.666
$ws.Range("N130").Value = -33306.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1993.3939
$ws.Range("I94").Value = 1906.3793
$ws.Range("J94").Value = 2624.25
$ws.Range("K94").Value = 1906.3793
$ws.Range("L94").Value = 2624.25
$ws.Range("M94").Value = -1455.3793
$ws.Range("N94").Value = -3526.25
$ws.Range("H105").Value = 68688.60000000001
$ws.Range("I105").Value = 112056
$ws.Range("K105").Value = 112056
$ws.Range("M105").Value = -110309
$ws.Range("H107").Value = 2134.8823
$ws.Range("I107").Value = 1644.9166
$ws.Range("J107").Value = 3310.8
$ws.Range("K107").Value = 1644.9166
$ws.Range("L107").Value = 3310.8
$ws.Range("M107").Value = 275.0834
$ws.Range("N107").Value = -7150.8
$ws.Range("H135").Value = 77775.55499999999
$ws.Range("J135").Value = 77775.55499999999
$ws.Range("L135").Value = 77775.55499999999
$ws.Range("N135").Value = -87915.55499999999
$ws.Range("H138").Value = 99752.28999999999
$ws.Range("J138").Value = 99752.28999999999
$ws.Range("L138").Value = 99752.28999999999
$ws.Range("N138").Value = -110032.29
$ws.Range("H140").Value = 43499.273
$ws.Range("J140").Value = 43499.273
$ws.Range("L140").Value = 43499.273
$ws.Range("N140").Value = -53859.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H138").Value = 94496
$ws.Range("J138").Value = 94496
$ws.Range("L138").Value = 94496
$ws.Range("N138").Value = -104776

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 7967.4
$ws.Range("J132").Value = 7967.4
$ws.Range("L132").Value = 71706.59999999999
$ws.Range("N132").Value = -76766.59999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2133.2
$ws.Range("I102").Value = 2041.5
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2041.5
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -419.5
$ws.Range("N102").Value = -5744
$ws.Range("H122").Value = 336998
$ws.Range("I122").Value = 479784.47
$ws.Range("J122").Value = 3829.5557
$ws.Range("K122").Value = 1439353.41
$ws.Range("L122").Value = 11488.6671
$ws.Range("M122").Value = -1436903.41
$ws.Range("N122").Value = -16388.6671
$ws.Range("H132").Value = 6578.3
$ws.Range("I132").Value = 6206.6
$ws.Range("J132").Value = 6950
$ws.Range("K132").Value = 18619.8
$ws.Range("L132").Value = 20850
$ws.Range("M132").Value = -16089.8
$ws.Range("N132").Value = -25910
$ws.Range("H135").Value = 42138.855
$ws.Range("J135").Value = 42138.855
$ws.Range("L135").Value = 42138.855
$ws.Range("N135").Value = -52278.855
$ws.Range("H140").Value = 90411.42999999999
$ws.Range("J140").Value = 90396.664
$ws.Range("L140").Value = 90396.664
$ws.Range("N140").Value = -100756.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4632116
$ws.Range("I40").Value = 2981.6667
$ws.Range("K40").Value = 2981.6667
$ws.Range("M40").Value = -2845.6667
$ws.Range("H82").Value = 2661
$ws.Range("J82").Value = 3492.5
$ws.Range("L82").Value = 3492.5
$ws.Range("N82").Value = -4214.5
$ws.Range("H85").Value = 2661
$ws.Range("J85").Value = 3492.5
$ws.Range("L85").Value = 3492.5
$ws.Range("N85").Value = -5988.5
$ws.Range("H136").Value = 2759.2
$ws.Range("I136").Value = 4050
$ws.Range("J136").Value = 2436.5
$ws.Range("K136").Value = 12150
$ws.Range("L136").Value = 7309.5
$ws.Range("M136").Value = -9600
$ws.Range("N136").Value = -12409.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 30812.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 30812.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 30812.5
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -31794.5
$ws.Range("H62").Value = 4216.3335
$ws.Range("I62").Value = 3824.75
$ws.Range("J62").Value = 4999.5
$ws.Range("K62").Value = 3824.75
$ws.Range("L62").Value = 4999.5
$ws.Range("M62").Value = -3200.75
$ws.Range("N62").Value = -6247.5
$ws.Range("H65").Value = 4216.3335
$ws.Range("I65").Value = 3824.75
$ws.Range("J65").Value = 4999.5
$ws.Range("K65").Value = 19123.75
$ws.Range("L65").Value = 24997.5
$ws.Range("M65").Value = -16003.75
$ws.Range("N65").Value = -31237.5
$ws.Range("H86").Value = 49999
$ws.Range("J86").Value = 49999
$ws.Range("L86").Value = 49999
$ws.Range("N86").Value = -52245
$ws.Range("H89").Value = 49999
$ws.Range("J89").Value = 49999
$ws.Range("L89").Value = 249995
$ws.Range("N89").Value = -261227
$ws.Range("H107").Value = 2300.6667
$ws.Range("I107").Value = 724.5
$ws.Range("J107").Value = 5453
$ws.Range("K107").Value = 2173.5
$ws.Range("L107").Value = 16359
$ws.Range("M107").Value = -253.5
$ws.Range("N107").Value = -20199
$ws.Range("H113").Value = 1108.1765
$ws.Range("I113").Value = 357
$ws.Range("K113").Value = 1071
$ws.Range("M113").Value = 1099
$ws.Range("H122").Value = 3149.8333
$ws.Range("I122").Value = 3099.75
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 9299.25
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -6849.25
$ws.Range("N122").Value = -14650
$ws.Range("H132").Value = 1297.3
$ws.Range("I132").Value = 903
$ws.Range("K132").Value = 2709
$ws.Range("M132").Value = -179
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
